$wb = $excel.ActiveWorkbook

# The data on sheets "展览" and "全部类型" is mirrored, and both need the
# same numeric updates to column F (the view/attendance count column).
$sheetNames = @("展览", "全部类型")

# Mapping of row -> new value for column F
$updates = @{
    9  = 5906
    13 = 1805
    14 = 1317
    16 = 505
    17 = 107
    18 = 5531
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
